# ------------------------------------------------------------------
# Refresh the NSE watch-list table: new symbols for the existing 20
# rows (B/C/E columns) and 31 additional rows (21-52) with new
# "support Zone" (column C) entries, extending the sheet from
# A1:F21 to A1:F52.
# ------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give the new index cells (A22:A52) the same look (bold, bordered,
# centered) as the existing A2:A21 index column by copying the
# formatting from an existing cell before filling in values.
$ws.Cells.Item(2, 1).Copy()
$ws.Range("A22:A52").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# entry layout: row index (0-based -> column A), B, C, D, E, F
$data = @(
    ,@(0, "NSE:DEEPINDS", "NSE:AIROLAM", "", "NSE:ATGL", "")
    ,@(1, "NSE:FDC", "NSE:ALMONDZ", "", "NSE:DIXON", "")
    ,@(2, "NSE:GOLDBEES", "NSE:BALAXI", "", "NSE:GRANULES", "")
    ,@(3, "NSE:GOLDETF", "NSE:BANKBARODA", "", "NSE:HDFCBANK", "")
    ,@(4, "NSE:IVZINGOLD", "NSE:BFINVEST", "", "NSE:HINDALCO", "")
    ,@(5, "NSE:SADBHIN", "NSE:BFUTILITIE", "", "NSE:INDUSINDBK", "")
    ,@(6, "", "NSE:BHAGYANGR", "", "NSE:IRFC", "")
    ,@(7, "", "NSE:BIRLAMONEY", "", "NSE:JINDALSTEL", "")
    ,@(8, "", "NSE:CAMPUS", "", "NSE:JSWSTEEL", "")
    ,@(9, "", "NSE:CLEDUCATE", "", "NSE:LICHSGFIN", "")
    ,@(10, "", "NSE:DCW", "", "NSE:MGL", "")
    ,@(11, "", "NSE:DIAMONDYD", "", "NSE:NCC", "")
    ,@(12, "", "NSE:DTIL", "", "", "")
    ,@(13, "", "NSE:DUCON", "", "", "")
    ,@(14, "", "NSE:EBBETF0430", "", "", "")
    ,@(15, "", "NSE:FACT", "", "", "")
    ,@(16, "", "NSE:FEDERALBNK", "", "", "")
    ,@(17, "", "NSE:FINCABLES", "", "", "")
    ,@(18, "", "NSE:FIVESTAR", "", "", "")
    ,@(19, "", "NSE:GENESYS", "", "", "")
    ,@(20, "", "NSE:GICRE", "", "", "")
    ,@(21, "", "NSE:GILT5YBEES", "", "", "")
    ,@(22, "", "NSE:GREENLAM", "", "", "")
    ,@(23, "", "NSE:GRINFRA", "", "", "")
    ,@(24, "", "NSE:GSLSU", "", "", "")
    ,@(25, "", "NSE:HARSHA", "", "", "")
    ,@(26, "", "NSE:HEG", "", "", "")
    ,@(27, "", "NSE:HEIDELBERG", "", "", "")
    ,@(28, "", "NSE:HINDCOMPOS", "", "", "")
    ,@(29, "", "NSE:HINDCOPPER", "", "", "")
    ,@(30, "", "NSE:IEX", "", "", "")
    ,@(31, "", "NSE:INDUSTOWER", "", "", "")
    ,@(32, "", "NSE:JASH", "", "", "")
    ,@(33, "", "NSE:KOTARISUG", "", "", "")
    ,@(34, "", "NSE:KREBSBIO", "", "", "")
    ,@(35, "", "NSE:MAGADSUGAR", "", "", "")
    ,@(36, "", "NSE:MANGCHEFER", "", "", "")
    ,@(37, "", "NSE:MANOMAY", "", "", "")
    ,@(38, "", "NSE:MHLXMIRU", "", "", "")
    ,@(39, "", "NSE:MINDACORP", "", "", "")
    ,@(40, "", "NSE:MUFIN", "", "", "")
    ,@(41, "", "NSE:MUKANDLTD", "", "", "")
    ,@(42, "", "NSE:MUTHOOTCAP", "", "", "")
    ,@(43, "", "NSE:NIITLTD", "", "", "")
    ,@(44, "", "NSE:ONWARDTEC", "", "", "")
    ,@(45, "", "NSE:PLASTIBLEN", "", "", "")
    ,@(46, "", "NSE:RAJSREESUG", "", "", "")
    ,@(47, "", "NSE:RATNAVEER", "", "", "")
    ,@(48, "", "NSE:REFEX", "", "", "")
    ,@(49, "", "NSE:RESPONIND", "", "", "")
    ,@(50, "", "NSE:SAGARDEEP", "", "", "")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $entry = $data[$i]

    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $ws.Cells.Item($row, 6).Value = $entry[5]
}
